$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values (M2:T2) ---
$ws.Range("M2").Value = 0.136464
$ws.Range("N2").Value = 0.409392
$ws.Range("O2").Value = 0.2657380151046518
$ws.Range("P2").Value = 0.2657380151046518
$ws.Range("Q2").Value = 0.0009276822719999999
$ws.Range("R2").Value = 0.008349140448
$ws.Range("S2").Value = 0.2657380151046518
$ws.Range("T2").Value = 0.2657380151046518

# --- Update existing row 3 values (O3,P3,S3,T3) ---
$ws.Range("O3").Value = 0.6222649188457632
$ws.Range("P3").Value = 0.6222649188457632
$ws.Range("S3").Value = 0.6222649188457632
$ws.Range("T3").Value = 0.6222649188457632

# --- Update existing row 4 values (K4:T4) ---
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01800166666666667
$ws.Range("N4").Value = 0.054005
$ws.Range("O4").Value = 0.03505486552186345
$ws.Range("P4").Value = 0.03505486552186345
$ws.Range("Q4").Value = 0.00012237533
$ws.Range("R4").Value = 0.00110137797
$ws.Range("S4").Value = 0.03505486552186345
$ws.Range("T4").Value = 0.03505486552186345

# --- Add new row 5 ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nxph1"
$ws.Range("C5").Value = "Nrxn2"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.006797999999999999
$ws.Range("H5").Value = 0.020394
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.039512
$ws.Range("N5").Value = 0.118536
$ws.Range("O5").Value = 0.0769422005277216
$ws.Range("P5").Value = 0.0769422005277216
$ws.Range("Q5").Value = 0.000268602576
$ws.Range("R5").Value = 0.002417423184
$ws.Range("S5").Value = 0.0769422005277216
$ws.Range("T5").Value = 0.0769422005277216
